$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the range to text format so numeric-looking strings (e.g. "1.00",
    # "0.750") are preserved exactly as typed instead of being parsed as
    # numbers, then restore the default "Normal" style so no stray
    # number-format style lingers on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "60.223.28"
$ws.Range("E2").Value = "  -5.21%  "
$ws.Range("E3").Value = "  -2.37%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue $ws.Range("D5") "565.77"
$ws.Range("E5").Value = "  -2.02%  "
Set-TextValue $ws.Range("D6") "131.88"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.341.63"
$ws.Range("E8").Value = "  -2.39%  "
Set-TextValue $ws.Range("D9") "0.474"
$ws.Range("E9").Value = "  -1.25%  "
Set-TextValue $ws.Range("D10") "7.45"
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "3.914.32"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "3.361.76"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("E16").Value = "  -2.74%  "
Set-TextValue $ws.Range("D17") "24.90"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "60.289.17"
$ws.Range("E18").Value = "  -5.18%  "
Set-TextValue $ws.Range("D19") "13.60"
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("E21").Value = "  -5.61%  "
Set-TextValue $ws.Range("D22") "355.87"
$ws.Range("E22").Value = "  -7.13%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.472.95"
$ws.Range("E24").Value = "  -2.52%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D25") "1.00"
$ws.Range("E25").Value = "  -0.02%  "
Set-TextValue $ws.Range("D26") "69.49"
$ws.Range("E26").Value = "  -5.78%  "
$ws.Range("E27").Value = "  +2.90%  "
Set-TextValue $ws.Range("D28") "1.65"
$ws.Range("E28").Value = "  +17.75%  "
Set-TextValue $ws.Range("D29") "7.52"
$ws.Range("E29").Value = "  +7.43%  "
Set-TextValue $ws.Range("D30") "0.999"
$ws.Range("E30").Value = "  -0.02%  "
Set-TextValue $ws.Range("D31") "8.00"
$ws.Range("E31").Value = "  +1.57%  "
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D35").Value = "3.372.45"
$ws.Range("E35").Value = "  -2.43%  "
Set-TextValue $ws.Range("D36") "22.95"
$ws.Range("E36").Value = "  +1.42%  "
Set-TextValue $ws.Range("D37") "5.38"
$ws.Range("E37").Value = "  +3.80%  "
$ws.Range("E38").Value = "  +2.56%  "
$ws.Range("E39").Value = "  +0.91%  "
Set-TextValue $ws.Range("D40") "158.84"
$ws.Range("E40").Value = "  -3.24%  "
$ws.Range("E41").Value = "  +1.59%  "
Set-TextValue $ws.Range("D42") "1.00"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  +2.28%  "
Set-TextValue $ws.Range("D44") "1.20"
$ws.Range("E44").Value = "  +9.42%  "
Set-TextValue $ws.Range("D45") "40.92"
$ws.Range("E45").Value = "  -0.82%  "
Set-TextValue $ws.Range("D46") "0.750"
$ws.Range("E46").Value = "  -4.27%  "
Set-TextValue $ws.Range("D47") "23.72"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("E48").Value = "  -0.06%  "
Set-TextValue $ws.Range("D49") "6.82"
$ws.Range("E49").Value = "  +1.88%  "
Set-TextValue $ws.Range("D50") "22.55"
$ws.Range("E50").Value = "  +11.61%  "
Set-TextValue $ws.Range("D51") "0.897"
$ws.Range("E51").Value = "  +2.04%  "
